$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Force the cell to hold the literal text even when it looks numeric
    # (e.g. "592.39", "1.00"), then drop back to the Normal style so we
    # don't leave a stray number-format behind on a cell that previously
    # had none.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "D2" "62.919.85"
$ws.Range("E2").Value = "  +2.46%  "

Set-TextCell "D3" "2.949.26"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextCell "D5" "592.39"
$ws.Range("E5").Value = "  -0.21%  "

Set-TextCell "D6" "147.23"
$ws.Range("E6").Value = "  +2.57%  "

Set-TextCell "D7" "1.00"
$ws.Range("E7").Value = "  -0.13%  "

Set-TextCell "D8" "2.949.68"
$ws.Range("E8").Value = "  +0.76%  "

$ws.Range("E9").Value = "  +1.15%  "

Set-TextCell "D10" "7.02"
$ws.Range("E10").Value = "  +1.19%  "

$ws.Range("E11").Value = "  +5.18%  "

Set-TextCell "D12" "0.436"
$ws.Range("E12").Value = "  -0.32%  "

Set-TextCell "D13" "0.0000233"
$ws.Range("E13").Value = "  +4.51%  "

Set-TextCell "D14" "32.44"
$ws.Range("E14").Value = "  -2.24%  "

$ws.Range("E15").Value = "  -1.34%  "

Set-TextCell "D16" "3.437.21"
$ws.Range("E16").Value = "  +0.56%  "

Set-TextCell "D17" "62.874.52"
$ws.Range("E17").Value = "  +2.44%  "

Set-TextCell "D18" "6.67"
$ws.Range("E18").Value = "  +0.72%  "

Set-TextCell "D19" "2.954.70"
$ws.Range("E19").Value = "  +0.81%  "

Set-TextCell "D20" "437.69"
$ws.Range("E20").Value = "  +1.06%  "

Set-TextCell "D21" "13.38"
$ws.Range("E21").Value = "  -1.12%  "

Set-TextCell "D22" "0.663"

Set-TextCell "D23" "7.00"
$ws.Range("E23").Value = "  -0.70%  "

Set-TextCell "D24" "11.24"
$ws.Range("E24").Value = "  +3.84%  "

Set-TextCell "D25" "80.65"
$ws.Range("E25").Value = "  -0.96%  "

$ws.Range("E26").Value = "  -2.16%  "

$ws.Range("E27").Value = "  +0.48%  "

$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("E29").Value = "  +1.64%  "

Set-TextCell "D30" "7.29"
$ws.Range("E30").Value = "  +6.04%  "

Set-TextCell "D31" "2.60"
$ws.Range("E31").Value = "  +0.46%  "

Set-TextCell "D32" "0.0₃0986"
$ws.Range("E32").Value = "  +13.09%  "

Set-TextCell "D33" "26.33"
$ws.Range("E33").Value = "  -1.41%  "

Set-TextCell "D34" "0.107"
$ws.Range("E34").Value = "  -0.78%  "

Set-TextCell "D35" "1.00"
$ws.Range("E35").Value = "  -0.06%  "

Set-TextCell "D36" "0.992"
$ws.Range("E36").Value = "  -1.60%  "

Set-TextCell "D37" "5.60"
$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("E38").Value = "  +1.10%  "

Set-TextCell "D39" "49.60"
$ws.Range("E39").Value = "  -0.27%  "

Set-TextCell "D40" "2.01"
$ws.Range("E40").Value = "  +1.46%  "

$ws.Range("E41").Value = "  -3.93%  "

$ws.Range("E42").Value = "  -0.82%  "

$ws.Range("E43").Value = "  -0.12%  "

Set-TextCell "D44" "39.47"
$ws.Range("E44").Value = "  -5.81%  "

$ws.Range("E45").Value = "  +1.62%  "

Set-TextCell "D46" "2.685.06"
$ws.Range("E46").Value = "  -0.35%  "

$ws.Range("E47").Value = "  -1.79%  "

Set-TextCell "D48" "357.48"
$ws.Range("E48").Value = "  -1.61%  "

$ws.Range("E50").Value = "  -0.58%  "

Set-TextCell "D51" "22.61"
$ws.Range("E51").Value = "  -3.84%  "
